# feature: PoC v1.2.0 - Versao final da Prova de Conceito para deployment
#
# The "machine" sheet lists calendars per machine in column H ("CalendarioId").
# ASHE1 (row 2) was still pointing at the old "CAL-PADRAO-5x8" calendar; the
# other machines already use "CAL-24x5". Point ASHE1 at "CAL-24x5" too so all
# machines share the same calendar ahead of deployment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "CAL-24x5"

# Leave the saved selection on the cell that was just edited.
$ws.Range("H2").Select()
